$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 1.75
$ws.Range("I4").Value = 4.33
$ws.Range("J4").Value = 2.3
$ws.Range("L4").Value = 4.33
$ws.Range("U4").Value = 1.57
$ws.Range("W4").Value = 10
$ws.Range("AA4").Value = 13
$ws.Range("AI4").Value = 23
$ws.Range("AM4").Value = 29
$ws.Range("AO4").Value = 9
$ws.Range("AY4").Value = 23

# Row 5
$ws.Range("G5").Value = 6.5
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 1.53
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 4.5
$ws.Range("Q5").Value = 1.65
$ws.Range("R5").Value = 2.2
$ws.Range("AD5").Value = 8

# Row 12
$ws.Range("G12").Value = 3.4
$ws.Range("H12").Value = 3.6
$ws.Range("I12").Value = 2.05
$ws.Range("J12").Value = 4
$ws.Range("L12").Value = 2.63
$ws.Range("O12").Value = 1.25
$ws.Range("P12").Value = 3.75
$ws.Range("Q12").Value = 1.88
$ws.Range("R12").Value = 1.98
$ws.Range("S12").Value = 1.36
$ws.Range("T12").Value = 3
$ws.Range("W12").Value = 11
$ws.Range("X12").Value = 19
$ws.Range("Y12").Value = 12
$ws.Range("Z12").Value = 41
$ws.Range("AA12").Value = 26
$ws.Range("AB12").Value = 34
$ws.Range("AD12").Value = 7
$ws.Range("AE12").Value = 15
$ws.Range("AH12").Value = 8
$ws.Range("AI12").Value = 10
$ws.Range("AK12").Value = 19
$ws.Range("AL12").Value = 15
$ws.Range("AN12").Value = 5.5
$ws.Range("AO12").Value = 19
$ws.Range("AQ12").Value = 67
$ws.Range("AT12").Value = 3
$ws.Range("AW12").Value = 4
$ws.Range("AX12").Value = 11
